# Swap the content (columns A,B,D,E,F,G,H,Q,R) between pairs/cycles of rows
# in the "Artfynd" sheet, per the commit's row-data shuffle.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that actually carry the changed data for these rows.
$cols = @(1, 2, 4, 5, 6, 7, 8, 17, 18)   # A, B, D, E, F, G, H, Q, R

function Get-RowValues($ws, $row, $cols) {
    $vals = @{}
    foreach ($c in $cols) {
        $vals[$c] = $ws.Cells.Item($row, $c).Value2
    }
    return $vals
}

function Set-RowValues($ws, $row, $cols, $vals) {
    foreach ($c in $cols) {
        $ws.Cells.Item($row, $c).Value2 = $vals[$c]
    }
}

# Rotation/swap groups: each entry lists rows whose contents cycle in order,
# i.e. row[0] gets row[1]'s old data, row[1] gets row[2]'s old data, ...,
# and the last row gets row[0]'s old data.
$groups = @(
    @(11, 13),
    @(12, 14),
    @(15, 17, 16),
    @(18, 19),
    @(40, 41)
)

foreach ($group in $groups) {
    $snapshot = @{}
    foreach ($r in $group) {
        $snapshot[$r] = Get-RowValues $ws $r $cols
    }
    for ($i = 0; $i -lt $group.Count; $i++) {
        $destRow = $group[$i]
        $srcRow = $group[($i + 1) % $group.Count]
        Set-RowValues $ws $destRow $cols $snapshot[$srcRow]
    }
}
